# Deployment config additions: add flaskSecret, flaskPort and mode rows
# beneath the existing appDbConStr / systemConstraintFolderPath rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new key/value pairs for rows 3-5, column A (keys) first,
# then column B (values) - matching how the workbook was originally edited.
$ws.Range("A3").Value = "flaskSecret"
$ws.Range("A4").Value = "flaskPort"
$ws.Range("A5").Value = "mode"

$ws.Range("B5").Value = "p"
$ws.Range("B4").Value = 80
$ws.Range("B3").Value = "sec"

# Leave the active selection on B5, as in the final workbook state.
$ws.Range("B5").Select()
